$wb = $excel.ActiveWorkbook

# The "想去人数" (want-to-go headcount) figures bumped for three events
# in both the "展览" and "全部类型" sheets (they mirror the same rows).
$sheetNames = @("展览", "全部类型")
$updates = @{
    "F4"  = 2865
    "F12" = 1273
    "F21" = 2859
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
